$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 9.983522426115931
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 30.34306516417429

# Row 3 values
$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.1496068669990043
$ws.Range("E3").Value = 2797.565817734744
$ws.Range("G3").Value = 2800.788059942304
